# Update PollsData with OpinionWay poll (11/18) - adds rows 111-113
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 111 ---
$ws.Cells.Item(111, 1).Value = 33
$ws.Cells.Item(111, 2).Value = 2021
$ws.Cells.Item(111, 3).Value = 12
$ws.Cells.Item(111, 4).Value = 11
$ws.Cells.Item(111, 5).Value = 16
$ws.Cells.Item(111, 6).Value = "opinionway"
$ws.Cells.Item(111, 7).Value = "online"
$ws.Cells.Item(111, 8).Value = "included"
$ws.Cells.Item(111, 9).Value = 1178
$ws.Cells.Item(111, 10).Value = 1
$ws.Cells.Item(111, 11).Value = 1
$ws.Cells.Item(111, 12).Value = 9
$ws.Cells.Item(111, 13).Value = 2
$ws.Cells.Item(111, 14).Value = 3
$ws.Cells.Item(111, 15).Value = 8
$ws.Cells.Item(111, 16).Value = 5
$ws.Cells.Item(111, 17).Value = 24
$ws.Cells.Item(111, 20).Value = 13
$ws.Cells.Item(111, 22).Value = 3
$ws.Cells.Item(111, 23).Value = 19
$ws.Cells.Item(111, 24).Value = 12

# --- Row 112 ---
$ws.Cells.Item(112, 1).Value = 33
$ws.Cells.Item(112, 2).Value = 2021
$ws.Cells.Item(112, 3).Value = 12
$ws.Cells.Item(112, 4).Value = 11
$ws.Cells.Item(112, 5).Value = 16
$ws.Cells.Item(112, 6).Value = "opinionway"
$ws.Cells.Item(112, 7).Value = "online"
$ws.Cells.Item(112, 8).Value = "included"
$ws.Cells.Item(112, 9).Value = 1178
$ws.Cells.Item(112, 10).Value = 1
$ws.Cells.Item(112, 11).Value = 1
$ws.Cells.Item(112, 12).Value = 9
$ws.Cells.Item(112, 13).Value = 2
$ws.Cells.Item(112, 14).Value = 3
$ws.Cells.Item(112, 15).Value = 8
$ws.Cells.Item(112, 16).Value = 6
$ws.Cells.Item(112, 17).Value = 25
$ws.Cells.Item(112, 18).Value = 11
$ws.Cells.Item(112, 22).Value = 3
$ws.Cells.Item(112, 23).Value = 19
$ws.Cells.Item(112, 24).Value = 12

# --- Row 113 ---
$ws.Cells.Item(113, 1).Value = 33
$ws.Cells.Item(113, 2).Value = 2021
$ws.Cells.Item(113, 3).Value = 12
$ws.Cells.Item(113, 4).Value = 11
$ws.Cells.Item(113, 5).Value = 16
$ws.Cells.Item(113, 6).Value = "opinionway"
$ws.Cells.Item(113, 7).Value = "online"
$ws.Cells.Item(113, 8).Value = "included"
$ws.Cells.Item(113, 9).Value = 1178
$ws.Cells.Item(113, 10).Value = 1
$ws.Cells.Item(113, 11).Value = 1
$ws.Cells.Item(113, 12).Value = 9
$ws.Cells.Item(113, 13).Value = 2
$ws.Cells.Item(113, 14).Value = 3
$ws.Cells.Item(113, 15).Value = 8
$ws.Cells.Item(113, 16).Value = 6
$ws.Cells.Item(113, 17).Value = 25
$ws.Cells.Item(113, 19).Value = 9
$ws.Cells.Item(113, 22).Value = 3
$ws.Cells.Item(113, 23).Value = 21
$ws.Cells.Item(113, 24).Value = 12

# --- View state: active cell / selection moves to J112, frozen pane scrolls down ---
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("J112").Select()
